$d = $word.ActiveDocument

# Locate the paragraph that contains the existing
# "{{ '<run>black</run>.png' | image() }}" placeholder so we can
# duplicate its run/formatting structure for the second image reference.
$srcParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "*black.png*") {
        $srcParaIndex = $i
    }
}

$srcPara = $d.Paragraphs.Item($srcParaIndex)
$srcRange = $d.Range($srcPara.Range.Start, $srcPara.Range.End)

# 1) Duplicate the original placeholder paragraph (with its exact run
#    formatting) twice right after itself - pasting (rather than
#    InsertParagraphAfter + setting .Text) keeps the paragraph-mark
#    formatting (pPr/rPr) empty/inherited, just like the source
#    paragraph, instead of materializing explicit run properties there.
$pasteTarget1 = $d.Range($srcPara.Range.End, $srcPara.Range.End)
$srcRange.Copy()
$pasteTarget1.Paste()

$pasteTarget2 = $d.Range($d.Paragraphs.Item($srcParaIndex + 1).Range.End, $d.Paragraphs.Item($srcParaIndex + 1).Range.End)
$pasteTarget2.Paste()

# 2) Turn the first duplicate into the new explanatory paragraph. Replace
#    the text only up to (but excluding) the trailing paragraph mark, so
#    the paragraph mark - and its formatting - stays untouched.
$descPara = $d.Paragraphs.Item($srcParaIndex + 1)
$descTextRange = $d.Range($descPara.Range.Start, $descPara.Range.End - 1)
$descTextRange.Text = "use the same image twice to ensure reading from the dummy works twice"

# 3) Fix up the quoting style on the second duplicated paragraph: the
#    original uses single smart quotes ('...'), the duplicate should use
#    double smart quotes ("..."paragraph).
$dupPara = $d.Paragraphs.Item($srcParaIndex + 2)

$leftDoubleQuote = [char]0x201C
$rightDoubleQuote = [char]0x201D

$chars = $dupPara.Range.Characters
for ($i = 1; $i -le $chars.Count; $i++) {
    $ch = $chars.Item($i)
    $code = 0
    if ($ch.Text.Length -gt 0) {
        $code = [int][char]($ch.Text)
    }
    if ($code -eq 0x2018) {
        $ch.Text = "$leftDoubleQuote"
    } elseif ($code -eq 0x2019) {
        $ch.Text = "$rightDoubleQuote"
    }
}
